$wb = $excel.ActiveWorkbook

# xlEdgeLeft=7, xlEdgeTop=8, xlEdgeBottom=9, xlEdgeRight=10
# xlContinuous=1 (thin line), xlLineStyleNone=-4142

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

$c1 = $ws1.Range("C1")
$d1 = $ws1.Range("D1")
$c1.ClearFormats()
$d1.ClearFormats()

$c1.Borders.Item(7).LineStyle = -4142
$c1.Borders.Item(10).LineStyle = -4142
$c1.Borders.Item(8).LineStyle = 1
$c1.Borders.Item(9).LineStyle = 1

$d1.Borders.Item(7).LineStyle = -4142
$d1.Borders.Item(10).LineStyle = 1
$d1.Borders.Item(8).LineStyle = 1
$d1.Borders.Item(9).LineStyle = 1

$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

$c1b = $ws2.Range("C1")
$d1b = $ws2.Range("D1")
$f1b = $ws2.Range("F1")
$g1b = $ws2.Range("G1")

$c1b.ClearFormats()
$d1b.ClearFormats()
$f1b.ClearFormats()
$g1b.ClearFormats()

$c1b.Borders.Item(7).LineStyle = -4142
$c1b.Borders.Item(10).LineStyle = -4142
$c1b.Borders.Item(8).LineStyle = 1
$c1b.Borders.Item(9).LineStyle = 1

$d1b.Borders.Item(7).LineStyle = -4142
$d1b.Borders.Item(10).LineStyle = 1
$d1b.Borders.Item(8).LineStyle = 1
$d1b.Borders.Item(9).LineStyle = 1

$f1b.Borders.Item(7).LineStyle = -4142
$f1b.Borders.Item(10).LineStyle = -4142
$f1b.Borders.Item(8).LineStyle = 1
$f1b.Borders.Item(9).LineStyle = 1

$g1b.Borders.Item(7).LineStyle = -4142
$g1b.Borders.Item(10).LineStyle = 1
$g1b.Borders.Item(8).LineStyle = 1
$g1b.Borders.Item(9).LineStyle = 1

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5 held an empty placeholder inline-string cell; drop it entirely.
$ws2.Range("G5").ClearContents()
